# Apply the updated crypto price / 1h-volume figures (and the swapped
# EnergySwap/Aptos rows) from the scheduled GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new text value. Kept as plain strings (matching the workbook's
# existing inlineStr/text cells) rather than letting Excel auto-convert
# numeric-looking text ("0.9981", "95.50", ...) into a real number, which
# would silently drop the original formatting (trailing zeros, etc.).
$updates = [ordered]@{
    D2 = '23.531.62'
    E2 = '  -0.77%  '
    D3 = '1.641.79'
    E3 = '  -0.71%  '
    D4 = '0.9981'
    E4 = '  -0.25%  '
    D5 = '0.9995'
    E5 = '  +0.05%  '
    D6 = '304.67'
    E6 = '  -0.57%  '
    D7 = '0.3802'
    E7 = '  +0.45%  '
    D8 = '51.74'
    E8 = '  -2.10%  '
    D9 = '0.3633'
    E9 = '  -1.29%  '
    D10 = '0.08190'
    E10 = '  +0.14%  '
    D11 = '1.238'
    E11 = '  -3.01%  '
    E12 = '  +0.31%  '
    D13 = '22.58'
    E13 = '  -2.58%  '
    D14 = '6.477'
    E14 = '  -3.83%  '
    D15 = '7.395'
    E15 = '  -0.39%  '
    D16 = '0.00001245'
    E16 = '  -2.40%  '
    D17 = '1.632.03'
    E17 = '  -1.15%  '
    D18 = '95.50'
    E18 = '  +0.08%  '
    D19 = '0.06950'
    E19 = '  +0.29%  '
    D20 = '6.588'
    E20 = '  -0.58%  '
    D21 = '17.54'
    E21 = '  -5.28%  '
    D22 = '0.9994'
    E22 = '  +0.10%  '
    D23 = '12.56'
    E23 = '  -3.64%  '
    D24 = '23.462.63'
    E24 = '  -1.11%  '
    D25 = '2.515'
    E25 = '  +3.51%  '
    D26 = '3.067'
    E26 = '  -5.70%  '
    D27 = '21.16'
    E27 = '  -1.46%  '
    D28 = '151.37'
    E28 = '  -0.58%  '
    D29 = '5.271'
    E29 = '  -0.58%  '
    D30 = '133.69'
    E30 = '  -2.67%  '
    D31 = '1.820.13'
    E31 = '  -0.95%  '
    E32 = '  -5.50%  '
    D33 = '6.659'
    E33 = '  -7.55%  '
    E34 = '  +7.23%  '
    D35 = '11.46'
    E35 = '  +3.38%  '
    D36 = '0.02766'
    E36 = '  -4.98%  '
    D37 = '0.2496'
    E37 = '  -3.88%  '
    D38 = '0.08788'
    E38 = '  -1.50%  '
    D39 = '0.07143'
    E39 = '  -3.27%  '
    D40 = '6.023'
    E40 = '  -6.05%  '
    D41 = '0.7069'
    E41 = '  -2.54%  '
    D42 = '1.345'
    E42 = '  -3.35%  '
    B43 = 'Aptos'
    C43 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    D43 = '12.19'
    E43 = '  -4.47%  '
    B44 = 'EnergySwap'
    C44 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    D44 = '15.82'
    E44 = '  -4.57%  '
    D45 = '0.6565'
    E45 = '  -2.06%  '
    D46 = '0.9996'
    E46 = '  +0.16%  '
    D47 = '2.284'
    E47 = '  -5.05%  '
    D48 = '3.965'
    E48 = '  -1.69%  '
    D49 = '0.07986'
    E49 = '  -1.31%  '
    D50 = '127.77'
    E50 = '  -1.39%  '
    D51 = '1.197'
    E51 = '  -3.11%  '
}

foreach ($addr in $updates.Keys) {
    $newValue = $updates[$addr]
    $range = $ws.Range($addr)
    if ($newValue -match "^[+-]?[0-9]*\.?[0-9]+$") {
        # Looks like a plain number to Excel's auto-detection -- prefix with
        # an apostrophe to force text entry, matching the target's inlineStr.
        $range.Value = "'" + $newValue
        # Re-applying the base style drops the quotePrefix flag Excel adds
        # for apostrophe-led text, so the cell format stays untouched.
        $range.Style = "Normal"
    } else {
        $range.Value = $newValue
    }
}
